$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.555.23"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.57"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -1.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.39"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.99"
$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07930"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.55"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.970"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.848.53"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.155"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.40"
$ws.Range("E17").Value = "  +2.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06682"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001036"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.550.27"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.406"
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.071.50"
$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.26"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.54"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.119"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.424"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.42"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9765"
$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09405"
$ws.Range("E33").Value = "  -0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.304"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.336"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06013"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.286"
$ws.Range("E39").Value = "  +3.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.182"
$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("E41").Value = "  -0.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5913"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1864"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.35"
$ws.Range("E44").Value = "  +2.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.242"
$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5589"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.913"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06693"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.23"
$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.049"
$ws.Range("E51").Value = "  -0.82%  "
